$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on cells whose new values look numeric,
# so Excel keeps them as literal text (matching the source data).
$textFormatCells = @(
    "D4",
    "D7",
    "D8",
    "D9",
    "D10",
    "D11",
    "D12",
    "D13",
    "D14",
    "D16",
    "D17",
    "D18",
    "D21",
    "D24",
    "D25",
    "D27",
    "D28",
    "D29",
    "D30",
    "D31",
    "D32",
    "D33",
    "D35",
    "D36",
    "D37",
    "D38",
    "D39",
    "D40",
    "D41",
    "D42",
    "D43",
    "D45",
    "D46",
    "D47",
    "D49",
    "D50",
    "D51"
)
foreach ($addr in $textFormatCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated values cell by cell.
$ws.Range('D2').Value = '27.248.75'
$ws.Range('E2').Value = '  +0.25%  '
$ws.Range('D3').Value = '1.907.34'
$ws.Range('E3').Value = '  +0.11%  '
$ws.Range('D4').Value = '1.001'
$ws.Range('E4').Value = '  +0.13%  '
$ws.Range('E5').Value = '  +0.44%  '
$ws.Range('E6').Value = '  +0.07%  '
$ws.Range('D7').Value = '0.5255'
$ws.Range('E7').Value = '  +0.32%  '
$ws.Range('D8').Value = '0.3818'
$ws.Range('E8').Value = '  +1.28%  '
$ws.Range('D9').Value = '0.07312'
$ws.Range('E9').Value = '  +0.85%  '
$ws.Range('D10').Value = '21.62'
$ws.Range('E10').Value = '  +2.05%  '
$ws.Range('D11').Value = '0.9059'
$ws.Range('E11').Value = '  +0.26%  '
$ws.Range('D12').Value = '0.08190'
$ws.Range('E12').Value = '  -3.90%  '
$ws.Range('D13').Value = '96.50'
$ws.Range('E13').Value = '  +0.06%  '
$ws.Range('D14').Value = '5.372'
$ws.Range('E14').Value = '  +1.45%  '
$ws.Range('D15').Value = '1.652.53'
$ws.Range('E15').Value = '  -13.61%  '
$ws.Range('D16').Value = '1.001'
$ws.Range('E16').Value = '  +0.06%  '
$ws.Range('D17').Value = '0.000008694'
$ws.Range('E17').Value = '  +0.62%  '
$ws.Range('D18').Value = '14.76'
$ws.Range('E18').Value = '  +1.27%  '
$ws.Range('E19').Value = '  +0.03%  '
$ws.Range('D20').Value = '27.279.41'
$ws.Range('E20').Value = '  +0.25%  '
$ws.Range('D21').Value = '5.126'
$ws.Range('E21').Value = '  +1.09%  '
$ws.Range('E23').Value = '  +1.17%  '
$ws.Range('D24').Value = '2.348'
$ws.Range('E24').Value = '  +2.17%  '
$ws.Range('D25').Value = '149.95'
$ws.Range('E25').Value = '  +1.61%  '
$ws.Range('E26').Value = '  -0.02%  '
$ws.Range('D27').Value = '1.736'
$ws.Range('E27').Value = '  -0.68%  '
$ws.Range('D28').Value = '117.02'
$ws.Range('E28').Value = '  +1.80%  '
$ws.Range('D29').Value = '4.846'
$ws.Range('E29').Value = '  +0.50%  '
$ws.Range('D30').Value = '4.870'
$ws.Range('E30').Value = '  -0.93%  '
$ws.Range('D31').Value = '0.09257'
$ws.Range('E31').Value = '  -0.36%  '
$ws.Range('D32').Value = '0.8241'
$ws.Range('E32').Value = '  +2.31%  '
$ws.Range('D33').Value = '0.05079'
$ws.Range('E33').Value = '  +0.43%  '
$ws.Range('E34').Value = '  -1.04%  '
$ws.Range('D35').Value = '2.984'
$ws.Range('E35').Value = '  +0.78%  '
$ws.Range('D36').Value = '2.748'
$ws.Range('E36').Value = '  +4.84%  '
$ws.Range('D37').Value = '3.360'
$ws.Range('E37').Value = '  -2.55%  '
$ws.Range('D38').Value = '0.5765'
$ws.Range('E38').Value = '  +1.00%  '
$ws.Range('D39').Value = '0.02004'
$ws.Range('E39').Value = '  +0.23%  '
$ws.Range('D40').Value = '1.083'
$ws.Range('E40').Value = '  +0.94%  '
$ws.Range('D41').Value = '9.079'
$ws.Range('E41').Value = '  -0.89%  '
$ws.Range('D42').Value = '6.599'
$ws.Range('E42').Value = '  -0.64%  '
$ws.Range('D43').Value = '116.91'
$ws.Range('E43').Value = '  +0.47%  '
$ws.Range('E44').Value = '  +0.27%  '
$ws.Range('D45').Value = '0.4919'
$ws.Range('E45').Value = '  +0.88%  '
$ws.Range('B46').Value = 'EnergySwap'
$ws.Range('C46').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D46').Value = '10.21'
$ws.Range('E46').Value = '  +0.69%  '
$ws.Range('B47').Value = 'PaxDollar'
$ws.Range('C47').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D47').Value = '1.001'
$ws.Range('E47').Value = '  +0.00%  '
$ws.Range('E48').Value = '  +1.76%  '
$ws.Range('D49').Value = '38.75'
$ws.Range('E49').Value = '  +3.13%  '
$ws.Range('B50').Value = 'Aave'
$ws.Range('C50').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D50').Value = '64.61'
$ws.Range('E50').Value = '  +0.48%  '
$ws.Range('B51').Value = 'Cronos'
$ws.Range('C51').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D51').Value = '0.06047'
$ws.Range('E51').Value = '  +1.76%  '
